$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 115.375
$ws.Range("I6").Value = 93.78570999999999
$ws.Range("K6").Value = 281.35713
$ws.Range("M6").Value = -169.35713

$ws.Range("H55").Value = 195.47368
$ws.Range("I55").Value = 128.07143
$ws.Range("J55").Value = 384.2
$ws.Range("K55").Value = 128.07143
$ws.Range("L55").Value = 384.2
$ws.Range("M55").Value = 85.92857000000001
$ws.Range("N55").Value = -812.2

$ws.Range("H62").Value = 4797.1113
$ws.Range("I62").Value = 4754.143
$ws.Range("J62").Value = 4947.5
$ws.Range("K62").Value = 4754.143
$ws.Range("L62").Value = 4947.5
$ws.Range("M62").Value = -4130.143
$ws.Range("N62").Value = -6195.5

$ws.Range("H65").Value = 4797.1113
$ws.Range("I65").Value = 4754.143
$ws.Range("J65").Value = 4947.5
$ws.Range("K65").Value = 23770.715
$ws.Range("L65").Value = 24737.5
$ws.Range("M65").Value = -20650.715
$ws.Range("N65").Value = -30977.5

$ws.Range("H129").Value = 458.7143
$ws.Range("I129").Value = 458.7143
$ws.Range("K129").Value = 1376.1429
$ws.Range("M129").Value = 3623.8571

$ws.Range("H138").Value = 58825920
$ws.Range("J138").Value = 4600
$ws.Range("L138").Value = 13800
$ws.Range("N138").Value = -24080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4644.5264
$ws.Range("I32").Value = 4666.467
$ws.Range("K32").Value = 4666.467
$ws.Range("M32").Value = -4379.467

$ws.Range("H45").Value = 14719.417
$ws.Range("I45").Value = 19762
$ws.Range("J45").Value = 7659.8
$ws.Range("K45").Value = 19762
$ws.Range("L45").Value = 7659.8
$ws.Range("M45").Value = -19385
$ws.Range("N45").Value = -8413.799999999999

$ws.Range("H55").Value = 22199.2
$ws.Range("J55").Value = 24999
$ws.Range("L55").Value = 24999
$ws.Range("N55").Value = -25629

$ws.Range("H74").Value = 1761.8572
$ws.Range("J74").Value = 3239.6667
$ws.Range("L74").Value = 3239.6667
$ws.Range("N74").Value = -4987.6667

$ws.Range("H77").Value = 1761.8572
$ws.Range("J77").Value = 3239.6667
$ws.Range("L77").Value = 16198.3335
$ws.Range("N77").Value = -24934.3335

$ws.Range("H122").Value = 17931.809
$ws.Range("I122").Value = 1768.4615
$ws.Range("J122").Value = 44197.25
$ws.Range("K122").Value = 5305.3845
$ws.Range("L122").Value = 132591.75
$ws.Range("M122").Value = -2855.3845
$ws.Range("N122").Value = -137491.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3148.6365
$ws.Range("J31").Value = 4998.6665
$ws.Range("L31").Value = 4998.6665
$ws.Range("N31").Value = -5588.6665

$ws.Range("H34").Value = 3148.6365
$ws.Range("J34").Value = 4998.6665
$ws.Range("L34").Value = 4998.6665
$ws.Range("N34").Value = -5402.6665

$ws.Range("H58").Value = 2063.8064
$ws.Range("I58").Value = 2006.826
$ws.Range("K58").Value = 2006.826
$ws.Range("M58").Value = -1803.826

$ws.Range("H110").Value = 149989.6
$ws.Range("J110").Value = 149989.6
$ws.Range("L110").Value = 149989.6
$ws.Range("N110").Value = -158169.6

$ws.Range("H134").Value = 3595.3845
$ws.Range("I134").Value = 3604.7
$ws.Range("K134").Value = 10814.1
$ws.Range("M134").Value = -8279.099999999999

$ws.Range("H136").Value = 2063.8064
$ws.Range("I136").Value = 2006.826
$ws.Range("K136").Value = 6020.478
$ws.Range("M136").Value = -3470.478

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 78.8125
$ws.Range("J12").Value = 106.09091
$ws.Range("L12").Value = 318.27273
$ws.Range("N12").Value = -664.2727299999999

$ws.Range("H56").Value = 19118.143
$ws.Range("I56").Value = 19118.143
$ws.Range("K56").Value = 19118.143
$ws.Range("M56").Value = -18588.143

$ws.Range("H68").Value = 939.1818
$ws.Range("I68").Value = 1097
$ws.Range("J68").Value = 663
$ws.Range("K68").Value = 3291
$ws.Range("L68").Value = 1989
$ws.Range("M68").Value = -2480
$ws.Range("N68").Value = -3611

$ws.Range("H71").Value = 939.1818
$ws.Range("I71").Value = 1097
$ws.Range("J71").Value = 663
$ws.Range("K71").Value = 9873
$ws.Range("L71").Value = 5967
$ws.Range("M71").Value = -5817
$ws.Range("N71").Value = -14079

$ws.Range("H92").Value = 714.2
$ws.Range("I92").Value = 757.6667
$ws.Range("J92").Value = 649
$ws.Range("K92").Value = 2273.0001
$ws.Range("L92").Value = 1947
$ws.Range("M92").Value = -1025.0001
$ws.Range("N92").Value = -4443

$ws.Range("H97").Value = 1082.2858
$ws.Range("J97").Value = 1304.5
$ws.Range("L97").Value = 3913.5
$ws.Range("N97").Value = -4905.5

$ws.Range("H109").Value = 620.3077
$ws.Range("I109").Value = 588.6667
$ws.Range("K109").Value = 1766.0001
$ws.Range("M109").Value = -726.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 99000
$ws.Range("I70").Value = 99000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 99000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -98730
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 99000
$ws.Range("I73").Value = 99000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 99000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -98064
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 36655.168
$ws.Range("I132").Value = 2374.6365
$ws.Range("K132").Value = 7123.9095
$ws.Range("M132").Value = -4593.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 34804
$ws.Range("I132").Value = 36182.914
$ws.Range("K132").Value = 108548.742
$ws.Range("M132").Value = -106018.742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12360
$ws.Range("J45").Value = 12229.091
$ws.Range("L45").Value = 12229.091
$ws.Range("N45").Value = -13211.091

$ws.Range("H122").Value = 6173
$ws.Range("I122").Value = 2517.7
$ws.Range("K122").Value = 7553.099999999999
$ws.Range("M122").Value = -5103.099999999999

$ws.Range("H132").Value = 3970.743
$ws.Range("I132").Value = 2978.1924
$ws.Range("K132").Value = 8934.5772
$ws.Range("M132").Value = -6404.5772

$ws.Range("H136").Value = 4283.516
$ws.Range("I136").Value = 4214.0713
$ws.Range("K136").Value = 12642.2139
$ws.Range("M136").Value = -10092.2139
